$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifting all existing data rows down by one)
$ws.Rows.Item(2).Insert()

# Fill in the new row 2 with the new accelerometer reading
$ws.Cells.Item(2, 1).Value = -2.385556173324585
$ws.Cells.Item(2, 2).Value = 2.812312602996826
$ws.Cells.Item(2, 3).Value = -0.1285117015242577

# Match formatting of the rest of the plain data rows (no special style)
$ws.Range("A2:C2").Style = "Normal"

# Remove the two oldest samples, which are now rows 22 and 23 after the insert
$ws.Range("A22:C23").EntireRow.Delete()
